# Weekly fruit/vegetable update: insert two new data rows (new market week)
# at the top of the data block (rows 938-939), pushing all existing rows
# down by two. Dimension grows from A1:R1014 to A1:R1016.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 938 (existing rows 938.. shift down to 940..)
$ws.Rows("938:939").Insert()

# New row 938: Acelga, Primera, Region Metropolitana, 2023-10-24 (serial 45223)
$ws.Cells.Item(938, 1).Value = 9
$ws.Cells.Item(938, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(938, 3).Value = "Metropolitana"
$ws.Cells.Item(938, 4).Value = 45223
$ws.Cells.Item(938, 5).Value = 13
$ws.Cells.Item(938, 6).Value = 100112009
$ws.Cells.Item(938, 7).Value = "Acelga"
$ws.Cells.Item(938, 8).Value = "Sin especificar"
$ws.Cells.Item(938, 9).Value = "Primera"
$ws.Cells.Item(938, 10).Value = 70
$ws.Cells.Item(938, 11).Value = 15000
$ws.Cells.Item(938, 12).Value = 15000
$ws.Cells.Item(938, 13).Value = 15000
$ws.Cells.Item(938, 14).Value = "$/docena de atados"
$ws.Cells.Item(938, 15).Value = "Región Metropolitana"
$ws.Cells.Item(938, 16).Value = 5000
$ws.Cells.Item(938, 17).Value = 3
$ws.Cells.Item(938, 18).Value = "Hortaliza"

# New row 939: Acelga, Segunda, Region Metropolitana, 2023-10-24 (serial 45223)
$ws.Cells.Item(939, 1).Value = 9
$ws.Cells.Item(939, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(939, 3).Value = "Metropolitana"
$ws.Cells.Item(939, 4).Value = 45223
$ws.Cells.Item(939, 5).Value = 13
$ws.Cells.Item(939, 6).Value = 100112009
$ws.Cells.Item(939, 7).Value = "Acelga"
$ws.Cells.Item(939, 8).Value = "Sin especificar"
$ws.Cells.Item(939, 9).Value = "Segunda"
$ws.Cells.Item(939, 10).Value = 52
$ws.Cells.Item(939, 11).Value = 12000
$ws.Cells.Item(939, 12).Value = 12000
$ws.Cells.Item(939, 13).Value = 12000
$ws.Cells.Item(939, 14).Value = "$/docena de atados"
$ws.Cells.Item(939, 15).Value = "Región Metropolitana"
$ws.Cells.Item(939, 16).Value = 4000
$ws.Cells.Item(939, 17).Value = 3
$ws.Cells.Item(939, 18).Value = "Hortaliza"
